$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AB1: new date header, styled like its row-1 neighbours (bold, thin box border, centered/top aligned)
$header = $ws.Range("AB1")
$header.Value = "13-10-2020"
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Font.Bold = $true

# Data cells AB2:AB36: plain numeric values (no special formatting), one new day of the series
$ws.Range("AB2").Value = 198
$ws.Range("AB3").Value = 43983
$ws.Range("AB4").Value = 2940
$ws.Range("AB5").Value = 28439
$ws.Range("AB6").Value = 10669
$ws.Range("AB7").Value = 1170
$ws.Range("AB8").Value = 27421
$ws.Range("AB9").Value = 102
$ws.Range("AB10").Value = 20535
$ws.Range("AB11").Value = 4465
$ws.Range("AB12").Value = 15414
$ws.Range("AB13").Value = 10401
$ws.Range("AB14").Value = 2637
$ws.Range("AB15").Value = 9992
$ws.Range("AB16").Value = 7776
$ws.Range("AB17").Value = 115795
$ws.Range("AB18").Value = 94473
$ws.Range("AB19").Value = 961
$ws.Range("AB20").Value = 14932
$ws.Range("AB21").Value = 212905
$ws.Range("AB22").Value = 2756
$ws.Range("AB23").Value = 2434
$ws.Range("AB24").Value = 156
$ws.Range("AB25").Value = 1409
$ws.Range("AB26").Value = 23430
$ws.Range("AB27").Value = 4617
$ws.Range("AB28").Value = 8576
$ws.Range("AB29").Value = 21671
$ws.Range("AB30").Value = 384
$ws.Range("AB31").Value = 43747
$ws.Range("AB32").Value = 24208
$ws.Range("AB33").Value = 3738
$ws.Range("AB34").Value = 6976
$ws.Range("AB35").Value = 38815
$ws.Range("AB36").Value = 30604

Write-Host "AB column populated"
